# iResearch (Dual Pie Chart) block — apply the commit's table/chart tweaks:
#   - table width 100% -> 99.86% (5000 -> 4993 fiftieths-of-a-percent)
#   - fixed table layout (so the shrink doesn't auto-expand back)
#   - tblGrid columns realigned to the cells' actual widths (5040/288/5040 dxa)
#   - chart row made taller (3177 -> 3600 twips) and vertically centered
#     so the charts stop touching the row's border lines
#   - chart default text Calibri 9pt -> 7pt ("All charts Calibri 7")

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- table width: w:tblW w:w="5000" w:type="pct" -> w:w="4993" ---
# Table.PreferredWidth (points) maps to the raw OOXML pct value / 20
# (5000/20 = 250.0 -> 4993/20 = 249.65).
$t.PreferredWidthType = 2   # wdPreferredWidthPercent
$t.PreferredWidth = 249.65

# --- add <w:tblLayout w:type="fixed"/> ---
$t.AllowAutoFit = $false

# --- tblGrid: 5072/239/5072 -> 5040/288/5040 ---
# Re-asserting each column's width (already matching every cell's tcW)
# makes Word regenerate tblGrid from the real cell widths instead of the
# stale grid, which is exactly the 5040/288/5040 dxa target.
$t.Columns.Item(1).Width = 252     # 5040 dxa
$t.Columns.Item(2).Width = 14.4    # 288 dxa
$t.Columns.Item(3).Width = 252     # 5040 dxa

# --- chart row: trHeight 3177 -> 3600 twips, cells vertically centered ---
$chartRow = $t.Rows.Item(2)
$chartRow.Height = 180             # 3600 twips = 180 pt

for ($i = 1; $i -le $chartRow.Cells.Count; $i++) {
    $cell = $chartRow.Cells.Item($i)
    $cell.VerticalAlignment = 1    # wdCellAlignVerticalCenter -> w:vAlign val="center"
}

# --- charts: default text run size 900 (9pt) -> 700 (7pt), Calibri ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ($shp.HasChart) {
        $chart = $shp.Chart
        $chart.ChartArea.Format.TextFrame2.TextRange.Font.Name = "Calibri"
        $chart.ChartArea.Format.TextFrame2.TextRange.Font.Size = 7
    }
}

Write-Output "iResearch dual pie chart block updated"
